$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.329.93'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '1.681.02'
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '''218.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").Value = '''0.5280'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.56%  '
$ws.Range("D7").Value = '''1.008'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '''0.2695'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.70%  '
$ws.Range("D9").Value = '''0.06467'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("D10").Value = '''21.97'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("D11").Value = '''0.07522'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.31%  '
$ws.Range("D12").Value = '1.686.24'
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").Value = '''4.516'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").Value = '''0.5790'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.20%  '
$ws.Range("D15").Value = '''0.000008508'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").Value = '''64.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.94%  '
$ws.Range("D17").Value = '26.342.50'
$ws.Range("E17").Value = '  +0.84%  '
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("E20").Value = '  +1.28%  '
$ws.Range("D21").Value = '''190.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.94%  '
$ws.Range("D22").Value = '''6.208'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = '''144.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("D25").Value = '''7.778'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.20%  '
$ws.Range("D26").Value = '''0.1256'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.69%  '
$ws.Range("D27").Value = '''15.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.44%  '
$ws.Range("D28").Value = '''0.06529'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.48%  '
$ws.Range("D29").Value = '''1.364'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.55%  '
$ws.Range("E30").Value = '  +0.86%  '
$ws.Range("D31").Value = '''3.594'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.02%  '
$ws.Range("D32").Value = '''3.590'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.52%  '
$ws.Range("D33").Value = '''1.660'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.98%  '
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("D35").Value = '''0.6220'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.81%  '
$ws.Range("D36").Value = '''2.405'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.67%  '
$ws.Range("D37").Value = '''2.739'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.00%  '
$ws.Range("D38").Value = '''6.296'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.71%  '
$ws.Range("D39").Value = '1.116.28'
$ws.Range("E39").Value = '  +4.04%  '
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("D41").Value = '''0.8755'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.98%  '
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("D43").Value = '''100.56'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").Value = '1.829.41'
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("D45").Value = '''0.00000000109'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.55%  '
$ws.Range("D46").Value = '''56.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.60%  '
$ws.Range("D47").Value = '''8.194'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.17%  '
$ws.Range("D48").Value = '''1.006'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("D49").Value = '''0.05268'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.18%  '
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("D51").Value = '''6.087'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.60%  '
